$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header values for new columns P (14) and Q (15) on row 1, matching style of O1
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

$hdr = $ws.Range("P1:Q1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.LineStyle = 1         # xlContinuous
$hdr.Borders.Weight = 2            # xlThin

# For rows 2-25: swap I<->K values, M becomes 2, O becomes 1, and add new P=2,Q=2 columns
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I
    $ws.Cells.Item($r, 11).Value = 1  # K
    $ws.Cells.Item($r, 13).Value = 2  # M
    $ws.Cells.Item($r, 15).Value = 1  # O
    $ws.Cells.Item($r, 16).Value = 2  # P
    $ws.Cells.Item($r, 17).Value = 2  # Q
}
